# Question6.xlsx edit:
#  - Add a new worksheet "Screener2" after "Filtered" that holds the
#    (previously combined) High/Low columns split into separate
#    "High" / "Low" columns plus the Name / Face Value columns -
#    i.e. the exception-handling / annotation pass described in the
#    commit message ("Attempted to handle the exception in question6
#    + Added some annotations").
#  - The new sheet becomes the active / selected sheet.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the last existing sheet ("Filtered")
# so the tab order becomes: Screener, Screener1, Filtered, Screener2.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Screener2"

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "High"
$ws.Range("C1").Value = "Low"
$ws.Range("D1").Value = "Face Value"

# Data row (first screened company, High/Low and Face Value split out
# into plain text values rather than the combined "High/Low" string)
$ws.Range("A2").Value = "Infosys Ltd"

# Force the numeric-looking values to be stored as text (matching the
# "1,733" / "1,305" / "5.00" literal strings) instead of being
# auto-converted to numbers with a thousands-separator number format.
$ws.Range("B2:D2").NumberFormat = "@"
$ws.Range("B2").Value = "1,733"
$ws.Range("C2").Value = "1,305"
$ws.Range("D2").Value = "5.00"
# Drop back to the workbook's default (Normal) style now that the
# text has been committed, so no stray number-format style lingers.
$ws.Range("B2:D2").Style = "Normal"

# Mirror the saved UI state: new sheet ends up active with H10 selected.
$ws.Activate() | Out-Null
$ws.Range("H10").Select() | Out-Null
